# Jenkins build details parameterize: update Regression (column C) flags
# on the Registration sheet, and move the active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

$regressionValues = @{
    2  = "NO"
    8  = "YES"
    9  = "YES"
    10 = "YES"
    11 = "YES"
    12 = "YES"
    14 = "YES"
    15 = "YES"
    16 = "YES"
    17 = "YES"
    18 = "YES"
    19 = "YES"
    20 = "YES"
    21 = "YES"
    24 = "YES"
    25 = "YES"
    26 = "YES"
    27 = "YES"
    28 = "YES"
    29 = "YES"
    30 = "YES"
    31 = "YES"
    32 = "YES"
    33 = "YES"
    34 = "YES"
    35 = "YES"
    36 = "YES"
    37 = "YES"
    40 = "YES"
    41 = "YES"
    42 = "YES"
    43 = "YES"
    44 = "YES"
    45 = "YES"
    46 = "YES"
    47 = "YES"
    48 = "YES"
    49 = "YES"
    50 = "YES"
    51 = "YES"
    52 = "YES"
    53 = "YES"
    56 = "YES"
    57 = "YES"
    58 = "YES"
    59 = "YES"
    60 = "YES"
    61 = "YES"
    62 = "YES"
    63 = "YES"
    64 = "YES"
    65 = "YES"
    66 = "YES"
    67 = "YES"
    68 = "YES"
    69 = "YES"
}

foreach ($row in $regressionValues.Keys) {
    $ws.Range("C$row").Value = $regressionValues[$row]
}

# Move the visible/selected window position like the author did while
# reviewing the freshly parameterized rows (scrolled to row 31, selected C58).
$ws.Activate()
$ws.Range("C58").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
